# Add two new bookings (SNOW-841107) as rows 7 and 8 on the active sheet,
# matching the existing row-5/row-6 "test" booking pattern.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 -------------------------------------------------------------
$ws.Range("A7").Value = "SNOW-841107"

# B7/K7 look like dates ("2026-02-19" / "2026-02-16") - format the cells as
# Text first so Excel stores the literal string instead of auto-converting
# it to a date serial number (matches how the other "test" rows are stored).
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2026-02-19"

$ws.Range("C7").Value = "test2"
$ws.Range("D7").Value = "test@test.com"
$ws.Range("E7").Value = "12134-34345"
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = "Family Ski Package"
$ws.Range("H7").Value = 32000
$ws.Range("I7").Value = 32000
$ws.Range("J7").Value = "Confirmed"

$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "2026-02-16"

# L7 is an empty "Special Requests" cell stored as an empty text value.
$ws.Range("L7").Formula = '=""'

# Row 8 -------------------------------------------------------------
$ws.Range("A8").Value = "SNOW-841107"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2026-02-20"

$ws.Range("C8").Value = "test2"
$ws.Range("D8").Value = "test@test.com"
$ws.Range("E8").Value = "12134-34345"
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = "Family Ski Package"
$ws.Range("H8").Value = 32000
$ws.Range("I8").Value = 32000
$ws.Range("J8").Value = "Confirmed"

$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "2026-02-16"

$ws.Range("L8").Formula = '=""'
